# Merge the three runs of the "Suspicious incident descriptions..." bullet
# on slide 4 into a single run with corrected punctuation/spacing:
#   "...can be red " + "flags.Evaluated" + " performance using:"
# becomes
#   "...can be red flags. Evaluated performance using"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# The bullet is the 7th paragraph of this placeholder's text.
$para = $tr.Paragraphs(7)

# Paragraph.Length includes the trailing paragraph mark, so Length-1 is the
# character count of the visible text. Re-assigning the whole paragraph's
# text (addressed via Characters, rather than via the Paragraph range
# itself) collapses it back down to a single run using the formatting of
# the paragraph's first run.
$fullText = $tr.Characters($para.Start, $para.Length - 1)
$fullText.Text = "Suspicious incident descriptions or vague reasons can be red flags. Evaluated performance using"
